$d = $word.ActiveDocument

# 1. Split the lone paragraph in two right after "Voici un magnifique mode
#    d'emploi." using Find/Replace so the existing (hidden) _GoBack
#    bookmark rides along with the tail of the text into the new second
#    paragraph, exactly as Word does when the user presses Enter there.
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$null = $find.Execute("d'emploi.", $false, $false, $false, $false, $false, $true, 1, $false, "d'emploi.^p", 2)

# 2. The _GoBack bookmark (created by Word, collapsed/empty) now sits at
#    the very start of the new second paragraph. Grow the new paragraph's
#    text on both sides of that bookmark so it ends up sandwiched between
#    two runs, matching "Il est malheureusement eno" + bookmark +
#    "cre très vide !".
$bm = $d.Bookmarks("_GoBack")
$bm.Range.InsertAfter("cre très vide !")

$bm = $d.Bookmarks("_GoBack")
$bm.Range.InsertBefore("Il est malheureusement eno")
